# Apply OCR-cleanup edits to the worksheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 14: strip stray pipe characters from header cells ---
$ws.Range("G14").Value = "数量"
$ws.Range("J14").Value = "希望納期"

# --- Row 15 ---
$ws.Range("I15").Value = "\1,200"
$ws.Range("K15").Value = "池田"
$ws.Range("L15").Value = "9/5ヤグチ精機様持込み"
$ws.Range("M15").Value = "間"
$ws.Range("N15").ClearContents()

# --- Row 16 ---
$ws.Range("I16").Value = "\1,200"
$ws.Range("K16").Value = "池田"
$ws.Range("L16").Value = "9/5ヤグチ精機様持込み"
$ws.Range("M16").Value = "に"
$ws.Range("N16").ClearContents()

# --- Row 17 ---
$ws.Range("I17").Value = "\1,200"
$ws.Range("K17").Value = "池田"
$ws.Range("L17").Value = "9/5ヤグチ精機様持込み"
$ws.Range("M17").Value = "拓"
$ws.Range("N17").ClearContents()

# --- Row 18 ---
$ws.Range("I18").Value = "\1,000"
$ws.Range("K18").Value = "池田"
$ws.Range("L18").Value = "9/5ヤグチ糖機様持込み"
$ws.Range("M18").ClearContents()

# --- Rows 20-22: drop stray pipe-only cells ---
$ws.Range("B20").ClearContents()
$ws.Range("C20").ClearContents()
$ws.Range("B21").ClearContents()
$ws.Range("C21").ClearContents()
$ws.Range("B22").ClearContents()
$ws.Range("C22").ClearContents()

# --- Row 25 ---
$ws.Range("F25").Value = "\368"
$ws.Range("G25").ClearContents()

# --- Row 26 ---
$ws.Range("F26").Value = "\4、968"
$ws.Range("G26").ClearContents()
$ws.Range("H26").ClearContents()
